# Fruta / hortaliza, semanal
#
# Inserts 5 new weekly price-report rows for "Uva" (grape) at the top of
# the existing data block (old row 821), pushing the previously-existing
# rows 821-878 down to 826-883. Also corrects three mis-ordered variety
# names that land at the new row numbers 843, 844 and 851 once the shift
# has happened.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Make room: insert 5 blank rows above the current row 821 --------
# This shifts old rows 821-878 down to 826-883 and extends the used range
# (dimension) from A1:T878 to A1:T883 automatically.
$ws.Rows("821:825").Insert()

# --- 2) Populate the 5 newly inserted rows with the new entries ---------
$newRows = @(
    @(6, 'Mercado Mayorista Lo Valledor de Santiago', 'Metropolitana', 44610, 13, 'Fruta', 100109, 'Uva', 100109001, 'Uva', 'Flame Seedless',    'Primera', 250, 10000, 10000, 10000, '$/bandeja 18 kilos', 'Región Metropolitana',    556, 18),
    @(6, 'Mercado Mayorista Lo Valledor de Santiago', 'Metropolitana', 44610, 13, 'Fruta', 100109, 'Uva', 100109001, 'Uva', 'Moscatel rosada',    'Primera', 250, 15000, 15000, 15000, '$/bandeja 18 kilos', 'Provincia de Limarí',     833, 18),
    @(6, 'Mercado Mayorista Lo Valledor de Santiago', 'Metropolitana', 44610, 13, 'Fruta', 100109, 'Uva', 100109001, 'Uva', 'Superior Seedless',  'Primera', 150, 10000, 10000, 10000, '$/bandeja 18 kilos', 'Región Metropolitana',    556, 18),
    @(6, 'Mercado Mayorista Lo Valledor de Santiago', 'Metropolitana', 44610, 13, 'Fruta', 100109, 'Uva', 100109001, 'Uva', 'Superior Seedless',  'Primera', 250, 10000, 10000, 10000, '$/bandeja 18 kilos', "Región de O'Higgins",     556, 18),
    @(6, 'Mercado Mayorista Lo Valledor de Santiago', 'Metropolitana', 44610, 13, 'Fruta', 100109, 'Uva', 100109001, 'Uva', 'Thompson seedless',  'Primera', 300, 12000, 12000, 12000, '$/bandeja 18 kilos', 'Región Metropolitana',    667, 18)
)

$startRow = 821
$r = $startRow
foreach ($row in $newRows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# --- 3) Fix three variety names (column K) that were swapped between ----
#        neighbouring rows; these land at the following row numbers
#        after the 5-row shift above.
$ws.Range("K843").Value = "Rosada pastilla"
$ws.Range("K844").Value = "Torontel blanca"
$ws.Range("K851").Value = "Red Globe"
